$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 826
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 72
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 42
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 1290
$ws.Range("L2").Value = 322
$ws.Range("M2").Value = 968
$ws.Range("N2").Value = 783
$ws.Range("O2").Value = 185
$ws.Range("P2").Value = 110
$ws.Range("Q2").Value = 67
$ws.Range("R2").Value = -43
$ws.Range("S2").Value = -55
$ws.Range("T2").Value = 10
$ws.Range("U2").Value = 57
$ws.Range("V2").Value = 155
$ws.Range("W2").Value = 6.08
$ws.Range("X2").Value = 6.41
$ws.Range("Y2").Value = 5.53
$ws.Range("Z2").Value = 4.14
$ws.Range("AA2").Value = 33.2
$ws.Range("AB2").Value = 613.05
$ws.Range("AC2").Value = 1917
$ws.Range("AD2").Value = 9.26
$ws.Range("AE2").Value = 39971
$ws.Range("AF2").Value = 0.44
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 1.97
$ws.Range("AI2").Value = 16.26
$ws.Range("AJ2").Value = 2200000

# Row 3
$ws.Range("D3").Value = 800
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 78
$ws.Range("H3").Value = 51
$ws.Range("I3").Value = 29
$ws.Range("J3").Value = 23
$ws.Range("K3").Value = 1333
$ws.Range("L3").Value = 313
$ws.Range("M3").Value = 1021
$ws.Range("N3").Value = 813
$ws.Range("O3").Value = 207
$ws.Range("P3").Value = 110
$ws.Range("Q3").Value = -3
$ws.Range("R3").Value = 36
$ws.Range("S3").Value = -17
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = -14
$ws.Range("V3").Value = 155
$ws.Range("W3").Value = 8.81
$ws.Range("X3").Value = 6.43
$ws.Range("Y3").Value = 3.59
$ws.Range("Z3").Value = 3.92
$ws.Range("AA3").Value = 30.63
$ws.Range("AB3").Value = 633.23
$ws.Range("AC3").Value = 1302
$ws.Range("AD3").Value = 16.2
$ws.Range("AE3").Value = 41513
$ws.Range("AF3").Value = 0.51
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 1.18
$ws.Range("AI3").Value = 17.1
$ws.Range("AJ3").Value = 2200000

# Row 4
$ws.Range("D4").Value = 854
$ws.Range("E4").Value = 104
$ws.Range("F4").Value = 104
$ws.Range("G4").Value = 141
$ws.Range("H4").Value = 104
$ws.Range("I4").Value = 70
$ws.Range("J4").Value = 34
$ws.Range("K4").Value = 1415
$ws.Range("L4").Value = 328
$ws.Range("M4").Value = 1087
$ws.Range("N4").Value = 862
$ws.Range("O4").Value = 226
$ws.Range("P4").Value = 110
$ws.Range("Q4").Value = 43
$ws.Range("R4").Value = 57
$ws.Range("S4").Value = -25
$ws.Range("T4").Value = 28
$ws.Range("U4").Value = 16
$ws.Range("V4").Value = 146
$ws.Range("W4").Value = 12.18
$ws.Range("X4").Value = 12.23
$ws.Range("Y4").Value = 8.41
$ws.Range("Z4").Value = 7.6
$ws.Range("AA4").Value = 30.14
$ws.Range("AB4").Value = 692.72
$ws.Range("AC4").Value = 3203
$ws.Range("AD4").Value = 7.13
$ws.Range("AE4").Value = 43985
$ws.Range("AF4").Value = 0.52
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 2.19
$ws.Range("AI4").Value = 13.9
$ws.Range("AJ4").Value = 2200000

# Row 5
$ws.Range("D5").Value = 741
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 72
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 37
$ws.Range("J5").Value = 18
$ws.Range("K5").Value = 1475
$ws.Range("L5").Value = 342
$ws.Range("M5").Value = 1133
$ws.Range("N5").Value = 892
$ws.Range("O5").Value = 241
$ws.Range("P5").Value = 110
$ws.Range("Q5").Value = 12
$ws.Range("R5").Value = -100
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = 41
$ws.Range("U5").Value = -30
$ws.Range("V5").Value = 158
$ws.Range("W5").Value = 6.72
$ws.Range("X5").Value = 7.44
$ws.Range("Y5").Value = 4.23
$ws.Range("Z5").Value = 3.82
$ws.Range("AA5").Value = 30.23
$ws.Range("AB5").Value = 717.77
$ws.Range("AC5").Value = 1685
$ws.Range("AD5").Value = 13.74
$ws.Range("AE5").Value = 45515
$ws.Range("AF5").Value = 0.51
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 1.08
$ws.Range("AI5").Value = 13.21
$ws.Range("AJ5").Value = 2200000

# Row 6
$ws.Range("D6").Value = 789
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = -9
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 374
$ws.Range("M6").Value = 1126
$ws.Range("N6").Value = 877
$ws.Range("P6").Value = 110
$ws.Range("Q6").Value = -24
$ws.Range("R6").Value = -4
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 18
$ws.Range("U6").Value = -43
$ws.Range("V6").Value = 183
$ws.Range("W6").Value = 2.5
$ws.Range("X6").Value = 0.58
$ws.Range("Y6").Value = -1.01
$ws.Range("Z6").Value = 0.31
$ws.Range("AA6").Value = 33.26
$ws.Range("AB6").Value = 747.4
$ws.Range("AC6").Value = -405
$ws.Range("AD6").Value = -46.04
$ws.Range("AE6").Value = 44744
$ws.Range("AF6").Value = 0.42
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 2200000

# Row 6: AG6 and AH6 no longer have data (removed)
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Rows 7-9: all data columns (D through AI) are removed, keep only A,B,C
$ws.Range("D7:AI9").ClearContents()
